# Append a new "2025-04-16" row (row 46) to every price sheet in the
# workbook, repeating the last known price (row 45, column B) for each
# sheet — matching the author's "Updated Argent prices" commit.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count
    $newRow = $lastRow + 1

    $lastPriceCell = $ws.Cells.Item($lastRow, 2)
    $dateCell = $ws.Cells.Item($newRow, 1)
    $priceCell = $ws.Cells.Item($newRow, 2)

    # Force text storage (matching the existing text-typed date/price
    # cells above) instead of letting Excel auto-coerce the date string
    # into a date serial or the numeric-looking price into a number.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = "2025-04-16"
    $priceCell.Value = $lastPriceCell.Text
}
